# Apply "Trade #192 closed / Trade #193 opened" update to the live trading
# results workbook.

$wb = $excel.ActiveWorkbook

# Helper: write a literal date-look-alike string (e.g. "2026-02-17") into a
# cell without Excel auto-converting it to a date serial number. We briefly
# force a Text number format so the string is stored as-is, then clear the
# format again so the cell ends up with the default (General) style - this
# matches how the source data is represented (plain string, no explicit
# style) while avoiding the automatic date parsing that a plain .Value
# assignment would trigger.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.46
$summary.Range("B4").Value = -0.54
$summary.Range("B5").Value = -0.06
$summary.Range("B6").Value = 192
$summary.Range("B8").Value = 80
$summary.Range("B9").Value = 41.67

# ---------------------------------------------------------------------
# Strategy Status sheet (volatility_scorer row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C12").Value = 99.27
$status.Range("D12").Value = 14
$status.Range("E12").Value = -0.73
$status.Range("F12").Value = -0.73
$status.Range("G12").Value = 35.71

# ---------------------------------------------------------------------
# Helper data for the two new trade rows
# ---------------------------------------------------------------------
# Trade #192 - volatility_scorer - CLOSED
$t192_A = 192
$t192_B = "2026-02-17"
$t192_C = "10:07:38"
$t192_D = "volatility_scorer"
$t192_E = "NEUTRAL"
$t192_F = 0.22
$t192_G = 0.1
$t192_H = "CLOSED"
$t192_I = -54.5455
$t192_J = -0.12
$t192_K = 99.27
$t192_L = 0
$t192_M = 0
$t192_N = 0.85
$t192_O = "Low vol market (score: inf) - ideal for market making"
$t192_P = "early_exit"
$t192_Q = 0.18

# Trade #193 - MarketMaking - OPEN
$t193_A = 193
$t193_B = "2026-02-17"
$t193_C = "10:07:38"
$t193_D = "MarketMaking"
$t193_E = "UP"
$t193_F = 0.78
$t193_H = "OPEN"
$t193_I = 0
$t193_J = 0
$t193_K = 100.1871991854615
$t193_L = 0
$t193_M = 0
$t193_N = 0.6
$t193_O = "Normal spread capture: 19600 bps"
$t193_Q = 0

# ---------------------------------------------------------------------
# All Trades sheet - append rows 193 (trade #192) and 194 (trade #193)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(193, 1).Value = $t192_A
Set-TextValue $allTrades.Cells.Item(193, 2) $t192_B
$allTrades.Cells.Item(193, 3).Value = $t192_C
$allTrades.Cells.Item(193, 4).Value = $t192_D
$allTrades.Cells.Item(193, 5).Value = $t192_E
$allTrades.Cells.Item(193, 6).Value = $t192_F
$allTrades.Cells.Item(193, 7).Value = $t192_G
$allTrades.Cells.Item(193, 8).Value = $t192_H
$allTrades.Cells.Item(193, 9).Value = $t192_I
$allTrades.Cells.Item(193, 10).Value = $t192_J
$allTrades.Cells.Item(193, 11).Value = $t192_K
$allTrades.Cells.Item(193, 12).Value = $t192_L
$allTrades.Cells.Item(193, 13).Value = $t192_M
$allTrades.Cells.Item(193, 14).Value = $t192_N
$allTrades.Cells.Item(193, 15).Value = $t192_O
$allTrades.Cells.Item(193, 16).Value = $t192_P
$allTrades.Cells.Item(193, 17).Value = $t192_Q

$allTrades.Cells.Item(194, 1).Value = $t193_A
Set-TextValue $allTrades.Cells.Item(194, 2) $t193_B
$allTrades.Cells.Item(194, 3).Value = $t193_C
$allTrades.Cells.Item(194, 4).Value = $t193_D
$allTrades.Cells.Item(194, 5).Value = $t193_E
$allTrades.Cells.Item(194, 6).Value = $t193_F
$allTrades.Cells.Item(194, 7).Value = ""
$allTrades.Cells.Item(194, 8).Value = $t193_H
$allTrades.Cells.Item(194, 9).Value = $t193_I
$allTrades.Cells.Item(194, 10).Value = $t193_J
$allTrades.Cells.Item(194, 11).Value = $t193_K
$allTrades.Cells.Item(194, 12).Value = $t193_L
$allTrades.Cells.Item(194, 13).Value = $t193_M
$allTrades.Cells.Item(194, 14).Value = $t193_N
$allTrades.Cells.Item(194, 15).Value = $t193_O
$allTrades.Cells.Item(194, 16).Value = ""
$allTrades.Cells.Item(194, 17).Value = $t193_Q

# ---------------------------------------------------------------------
# volatility_scorer sheet - append row 15 (trade #192)
# ---------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")

$volScorer.Cells.Item(15, 1).Value = $t192_A
Set-TextValue $volScorer.Cells.Item(15, 2) $t192_B
$volScorer.Cells.Item(15, 3).Value = $t192_C
$volScorer.Cells.Item(15, 4).Value = $t192_D
$volScorer.Cells.Item(15, 5).Value = $t192_E
$volScorer.Cells.Item(15, 6).Value = $t192_F
$volScorer.Cells.Item(15, 7).Value = $t192_G
$volScorer.Cells.Item(15, 8).Value = $t192_H
$volScorer.Cells.Item(15, 9).Value = $t192_I
$volScorer.Cells.Item(15, 10).Value = $t192_J
$volScorer.Cells.Item(15, 11).Value = $t192_K
$volScorer.Cells.Item(15, 12).Value = $t192_L
$volScorer.Cells.Item(15, 13).Value = $t192_M
$volScorer.Cells.Item(15, 14).Value = $t192_N
$volScorer.Cells.Item(15, 15).Value = $t192_O
$volScorer.Cells.Item(15, 16).Value = $t192_P
$volScorer.Cells.Item(15, 17).Value = $t192_Q

# ---------------------------------------------------------------------
# MarketMaking sheet - append row 180 (trade #193)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Cells.Item(180, 1).Value = $t193_A
Set-TextValue $marketMaking.Cells.Item(180, 2) $t193_B
$marketMaking.Cells.Item(180, 3).Value = $t193_C
$marketMaking.Cells.Item(180, 4).Value = $t193_D
$marketMaking.Cells.Item(180, 5).Value = $t193_E
$marketMaking.Cells.Item(180, 6).Value = $t193_F
$marketMaking.Cells.Item(180, 7).Value = ""
$marketMaking.Cells.Item(180, 8).Value = $t193_H
$marketMaking.Cells.Item(180, 9).Value = $t193_I
$marketMaking.Cells.Item(180, 10).Value = $t193_J
$marketMaking.Cells.Item(180, 11).Value = $t193_K
$marketMaking.Cells.Item(180, 12).Value = $t193_L
$marketMaking.Cells.Item(180, 13).Value = $t193_M
$marketMaking.Cells.Item(180, 14).Value = $t193_N
$marketMaking.Cells.Item(180, 15).Value = $t193_O
$marketMaking.Cells.Item(180, 16).Value = ""
$marketMaking.Cells.Item(180, 17).Value = $t193_Q
